# Auto-generated edit script: update calculated profit columns (H-N) across all profession sheets
# per scheduled-runner refresh of Sheets/Maduin_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1886.125
$ws.Range("I15").Value = 1886.125
$ws.Range("K15").Value = 5658.375
$ws.Range("M15").Value = -5489.375
$ws.Range("H98").Value = 2460.625
$ws.Range("I98").Value = 2460.625
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2460.625
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -962.625
$ws.Range("N98").ClearContents()
$ws.Range("H103").Value = 5582.75
$ws.Range("J103").Value = 6110.3335
$ws.Range("L103").Value = 18331.0005
$ws.Range("N103").Value = -19503.0005
$ws.Range("H113").Value = 4648
$ws.Range("I113").Value = 4500
$ws.Range("K113").Value = 4500
$ws.Range("M113").Value = -1246
$ws.Range("H122").Value = 2460.625
$ws.Range("I122").Value = 2460.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7381.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4931.875
$ws.Range("N122").ClearContents()
$ws.Range("H137").Value = 2648.1667
$ws.Range("I137").Value = 1966.3334
$ws.Range("J137").Value = 3330
$ws.Range("K137").Value = 5899.0002
$ws.Range("L137").Value = 9990
$ws.Range("M137").Value = -3349.0002
$ws.Range("N137").Value = -15090
$ws.Range("H138").Value = 14181.621
$ws.Range("J138").Value = 14181.621
$ws.Range("L138").Value = 42544.863
$ws.Range("N138").Value = -52824.863
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11992.667
$ws.Range("I32").Value = 11342.3
$ws.Range("K32").Value = 11342.3
$ws.Range("M32").Value = -11055.3
$ws.Range("H45").Value = 2697.3125
$ws.Range("I45").Value = 1519.625
$ws.Range("K45").Value = 1519.625
$ws.Range("M45").Value = -1142.625
$ws.Range("H74").Value = 3874.75
$ws.Range("I74").Value = 3874.75
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 3874.75
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3000.75
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 3874.75
$ws.Range("I77").Value = 3874.75
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 19373.75
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -15005.75
$ws.Range("N77").ClearContents()
$ws.Range("H110").Value = 3355.5881
$ws.Range("I110").Value = 2672.1667
$ws.Range("K110").Value = 2672.1667
$ws.Range("M110").Value = -627.1667000000002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 464.33334
$ws.Range("I5").Value = 197
$ws.Range("J5").Value = 731.6667
$ws.Range("K5").Value = 197
$ws.Range("L5").Value = 731.6667
$ws.Range("M5").Value = -84
$ws.Range("N5").Value = -957.6667
$ws.Range("H36").Value = 1604.375
$ws.Range("I36").Value = 1119.2858
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 1119.2858
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = -585.2858000000001
$ws.Range("N36").Value = -6068
$ws.Range("H105").Value = 5050
$ws.Range("I105").Value = 5050
$ws.Range("K105").Value = 5050
$ws.Range("M105").Value = -3303
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1588.5883
$ws.Range("I7").Value = 1345.625
$ws.Range("J7").Value = 1804.5555
$ws.Range("K7").Value = 1345.625
$ws.Range("L7").Value = 1804.5555
$ws.Range("M7").Value = -1232.625
$ws.Range("N7").Value = -2030.5555
$ws.Range("H16").Value = 540.75
$ws.Range("I16").Value = 540.75
$ws.Range("K16").Value = 540.75
$ws.Range("M16").Value = -253.75
$ws.Range("H31").Value = 2835.0833
$ws.Range("I31").Value = 2502.3
$ws.Range("K31").Value = 2502.3
$ws.Range("M31").Value = -2207.3
$ws.Range("H34").Value = 2835.0833
$ws.Range("I34").Value = 2502.3
$ws.Range("K34").Value = 2502.3
$ws.Range("M34").Value = -2300.3
$ws.Range("H69").Value = 12235.125
$ws.Range("I69").Value = 12235.125
$ws.Range("K69").Value = 12235.125
$ws.Range("M69").Value = -11486.125
$ws.Range("H72").Value = 12235.125
$ws.Range("I72").Value = 12235.125
$ws.Range("K72").Value = 36705.375
$ws.Range("M72").Value = -32961.375
$ws.Range("H113").Value = 540.75
$ws.Range("I113").Value = 540.75
$ws.Range("K113").Value = 540.75
$ws.Range("M113").Value = 1629.25
$ws.Range("H132").Value = 3959.7144
$ws.Range("I132").Value = 3110.4614
$ws.Range("K132").Value = 9331.3842
$ws.Range("M132").Value = -6801.3842
$ws.Range("H141").Value = 587489.9
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 587489.9
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 587489.9
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -597849.9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 39.57143
$ws.Range("I2").Value = 13.777778
$ws.Range("K2").Value = 82.666668
$ws.Range("M2").Value = 30.333332
$ws.Range("H12").Value = 4520
$ws.Range("J12").Value = 4520
$ws.Range("L12").Value = 13560
$ws.Range("N12").Value = -13906
$ws.Range("H40").Value = 30.666666
$ws.Range("I40").Value = 45
$ws.Range("K40").Value = 180
$ws.Range("M40").Value = -111
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1742.2
$ws.Range("I102").Value = 904
$ws.Range("J102").Value = 2999.5
$ws.Range("K102").Value = 904
$ws.Range("L102").Value = 2999.5
$ws.Range("M102").Value = 718
$ws.Range("N102").Value = -6243.5
$ws.Range("H132").Value = 3634.625
$ws.Range("I132").Value = 3101.2222
$ws.Range("K132").Value = 9303.6666
$ws.Range("M132").Value = -6773.6666
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3147.25
$ws.Range("J7").Value = 3274.5
$ws.Range("L7").Value = 3274.5
$ws.Range("N7").Value = -3498.5
$ws.Range("H46").Value = 1319.9
$ws.Range("I46").Value = 1139.6
$ws.Range("J46").Value = 1500.2
$ws.Range("K46").Value = 1139.6
$ws.Range("L46").Value = 1500.2
$ws.Range("M46").Value = -951.5999999999999
$ws.Range("N46").Value = -1876.2
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H122").Value = 7267.636
$ws.Range("I122").Value = 7594.5
$ws.Range("K122").Value = 22783.5
$ws.Range("M122").Value = -20333.5
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 3147.25
$ws.Range("J126").Value = 3274.5
$ws.Range("L126").Value = 9823.5
$ws.Range("N126").Value = -14763.5
$ws.Range("H132").Value = 8334
$ws.Range("I132").Value = 3003
$ws.Range("J132").Value = 10999.5
$ws.Range("K132").Value = 9009
$ws.Range("L132").Value = 32998.5
$ws.Range("M132").Value = -6479
$ws.Range("N132").Value = -38058.5
$ws.Range("H136").Value = 5000.625
$ws.Range("I136").Value = 5166.6665
$ws.Range("J136").Value = 4502.5
$ws.Range("K136").Value = 15499.9995
$ws.Range("L136").Value = 13507.5
$ws.Range("M136").Value = -12949.9995
$ws.Range("N136").Value = -18607.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 44500
$ws.Range("J33").Value = 44500
$ws.Range("L33").Value = 44500
$ws.Range("N33").Value = -45000
$ws.Range("H36").Value = 44500
$ws.Range("J36").Value = 44500
$ws.Range("L36").Value = 44500
$ws.Range("N36").Value = -45000
$ws.Range("H100").Value = 6973555
$ws.Range("I100").Value = 8715569
$ws.Range("J100").Value = 5499.5
$ws.Range("K100").Value = 17431138
$ws.Range("L100").Value = 10999
$ws.Range("M100").Value = -17430597
$ws.Range("N100").Value = -12081
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 2014.9
$ws.Range("I126").Value = 1518
$ws.Range("J126").Value = 4002.5
$ws.Range("K126").Value = 4554
$ws.Range("L126").Value = 12007.5
$ws.Range("M126").Value = -2084
$ws.Range("N126").Value = -16947.5
$ws.Range("H132").Value = 2189.7222
$ws.Range("I132").Value = 1507.9286
$ws.Range("J132").Value = 4576
$ws.Range("K132").Value = 4523.7858
$ws.Range("L132").Value = 13728
$ws.Range("M132").Value = -1993.7858
$ws.Range("N132").Value = -18788
$ws.Range("H136").Value = 40856.77
$ws.Range("I136").Value = 46921.637
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 140764.911
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -138214.911
$ws.Range("N136").Value = -27600
